# Updates cryptos list: refreshed Price (D) and Volume(1h) (E) columns
# Values that look numeric are apostrophe-prefixed so Excel keeps them as
# literal text (matching the source data, e.g. "43.00" not 43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.558.10"
$ws.Range("E2").Value = "  +1.57%  "

$ws.Range("D3").Value = "1.941.49"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'243.33"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("E6").Value = "  +1.50%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "'57.17"
$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("D9").Value = "'0.363"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "'0.0801"
$ws.Range("E10").Value = "  -1.98%  "

$ws.Range("E11").Value = "  -0.28%  "

$ws.Range("E12").Value = "  +2.94%  "

$ws.Range("D13").Value = "2.226.65"
$ws.Range("E13").Value = "  +0.33%  "

$ws.Range("E14").Value = "  -1.78%  "

$ws.Range("D15").Value = "'13.29"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").Value = "1.942.09"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").Value = "36.467.99"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").Value = "'69.22"
$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").Value = "0.0₃0853"
$ws.Range("E20").Value = "  -0.85%  "

$ws.Range("D21").Value = "'227.68"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").Value = "'2.39"
$ws.Range("E24").Value = "  -2.87%  "

$ws.Range("E25").Value = "  +1.34%  "

$ws.Range("E26").Value = "  -1.20%  "

$ws.Range("D27").Value = "'159.41"
$ws.Range("E27").Value = "  -2.04%  "

$ws.Range("D28").Value = "'0.135"
$ws.Range("E28").Value = "  +17.03%  "

$ws.Range("D29").Value = "'19.21"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").Value = "  +1.17%  "

$ws.Range("E31").Value = "  -3.66%  "

$ws.Range("D32").Value = "'4.62"
$ws.Range("E32").Value = "  -0.67%  "

$ws.Range("D33").Value = "'0.0616"
$ws.Range("E33").Value = "  -0.96%  "

$ws.Range("D34").Value = "'4.18"
$ws.Range("E34").Value = "  -1.59%  "

$ws.Range("D35").Value = "'6.14"
$ws.Range("E35").Value = "  +1.91%  "

$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("E38").Value = "  +3.10%  "

$ws.Range("D39").Value = "'3.27"
$ws.Range("E39").Value = "  +15.69%  "

$ws.Range("D40").Value = "'0.0984"
$ws.Range("E40").Value = "  +2.61%  "

$ws.Range("D41").Value = "'2.91"
$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("E44").Value = "  +2.14%  "

$ws.Range("D45").Value = "1.343.65"
$ws.Range("E45").Value = "  +1.39%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("E47").Value = "  -0.99%  "

$ws.Range("D48").Value = "'7.12"
$ws.Range("E48").Value = "  -2.07%  "

$ws.Range("E49").Value = "  +0.33%  "

$ws.Range("D50").Value = "2.118.35"
$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("D51").Value = "'43.00"
$ws.Range("E51").Value = "  -5.19%  "
